$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set all cell values to match the target data (new plotting dataset from
# Dr. Creuziger Matlab analysis: adds OffsetF/OffsetA/RD Single/TD Single rows
# and a new 1Pair-B column, and recomputes all averaged intensity ratios).

$ws.Range('B1').Value = 0
$ws.Range('C1').Value = 1
$ws.Range('D1').Value = 2
$ws.Range('E1').Value = 3
$ws.Range('F1').Value = 4
$ws.Range('G1').Value = 5
$ws.Range('H1').Value = 6
$ws.Range('I1').Value = 7
$ws.Range('J1').Value = 8
$ws.Range('K1').Value = 9
$ws.Range('L1').Value = 10
$ws.Range('M1').Value = 11
$ws.Range('N1').Value = 12
$ws.Range('O1').Value = 13
$ws.Range('P1').Value = 14
$ws.Range('Q1').Value = 15
$ws.Range('R1').Value = 16
$ws.Range('S1').Value = 17
$ws.Range('T1').Value = 18
$ws.Range('A2').Value = 0
$ws.Range('B2').Value = 'HKL'
$ws.Range('C2').Value = '[1, 1, 0]'
$ws.Range('D2').Value = '[2, 0, 0]'
$ws.Range('E2').Value = '[2, 1, 1]'
$ws.Range('F2').Value = '[2, 2, 0]'
$ws.Range('G2').Value = '[3, 1, 0]'
$ws.Range('H2').Value = '[2, 2, 2]'
$ws.Range('I2').Value = '[3, 2, 1]'
$ws.Range('J2').Value = '[4, 0, 0]'
$ws.Range('K2').Value = '1Pair-A'
$ws.Range('L2').Value = '1Pair-B'
$ws.Range('M2').Value = '2Pairs-A'
$ws.Range('N2').Value = '2Pairs-B'
$ws.Range('O2').Value = '3Pairs-A'
$ws.Range('P2').Value = '3Pairs-B'
$ws.Range('Q2').Value = '3Pairs-C'
$ws.Range('R2').Value = '4Pairs'
$ws.Range('S2').Value = '5A4F'
$ws.Range('T2').Value = 'MaxUnique'
$ws.Range('A3').Value = 1
$ws.Range('B3').Value = 'Equal Angle'
$ws.Range('C3').Value = 0.9896613832853026
$ws.Range('D3').Value = 0.8257276657060518
$ws.Range('E3').Value = 1.060648414985591
$ws.Range('F3').Value = 0.9896613832853026
$ws.Range('G3').Value = 0.8822622478386167
$ws.Range('H3').Value = 1.189157060518732
$ws.Range('I3').Value = 1.042946685878962
$ws.Range('J3').Value = 0.8257276657060518
$ws.Range('K3').Value = 0.9896613832853026
$ws.Range('L3').Value = 1.060648414985591
$ws.Range('M3').Value = 0.9431880403458213
$ws.Range('N3').Value = 0.9431880403458213
$ws.Range('O3').Value = 0.9228794428434197
$ws.Range('P3').Value = 0.9586791546589817
$ws.Range('Q3').Value = 0.9586791546589817
$ws.Range('R3').Value = 0.9664247118155619
$ws.Range('S3').Value = 0.9664247118155619
$ws.Range('T3').Value = 0.9984005763688759
$ws.Range('A4').Value = 2
$ws.Range('B4').Value = 'CLR'
$ws.Range('C4').Value = 1.006458699853815
$ws.Range('D4').Value = 0.9702045769904586
$ws.Range('E4').Value = 0.9977272791950519
$ws.Range('F4').Value = 1.006458699853815
$ws.Range('G4').Value = 0.9801627344677278
$ws.Range('H4').Value = 1.001297630404104
$ws.Range('I4').Value = 1.000221549211721
$ws.Range('J4').Value = 0.9702045769904586
$ws.Range('K4').Value = 1.006458699853815
$ws.Range('L4').Value = 0.9977272791950519
$ws.Range('M4').Value = 0.9839659280927553
$ws.Range('N4').Value = 0.9839659280927553
$ws.Range('O4').Value = 0.9826981968844128
$ws.Range('P4').Value = 0.9914635186797751
$ws.Range('Q4').Value = 0.9914635186797751
$ws.Range('R4').Value = 0.995212313973285
$ws.Range('S4').Value = 0.995212313973285
$ws.Range('T4').Value = 0.9926787450204797
$ws.Range('A5').Value = 3
$ws.Range('B5').Value = 'BT8Hex'
$ws.Range('C5').Value = 1.015325854287473
$ws.Range('D5').Value = 0.9520672608689356
$ws.Range('E5').Value = 0.9999041670986337
$ws.Range('F5').Value = 1.015325854287473
$ws.Range('G5').Value = 0.9716656593299965
$ws.Range('H5').Value = 1.006017366443895
$ws.Range('I5').Value = 1.004408743195561
$ws.Range('J5').Value = 0.9520672608689356
$ws.Range('K5').Value = 1.015325854287473
$ws.Range('L5').Value = 0.9999041670986337
$ws.Range('M5').Value = 0.9759857139837846
$ws.Range('N5').Value = 0.9759857139837846
$ws.Range('O5').Value = 0.9745456957658553
$ws.Range('P5').Value = 0.9890990940850143
$ws.Range('Q5').Value = 0.9890990940850143
$ws.Range('R5').Value = 0.9956557841356291
$ws.Range('S5').Value = 0.9956557841356291
$ws.Range('T5').Value = 0.9915648418707494
$ws.Range('A6').Value = 4
$ws.Range('B6').Value = 'Spiral'
$ws.Range('C6').Value = 0.9919530836898475
$ws.Range('D6').Value = 0.9927867998959067
$ws.Range('E6').Value = 0.9967543493440122
$ws.Range('F6').Value = 0.9919530836898475
$ws.Range('G6').Value = 0.9902076254589198
$ws.Range('H6').Value = 1.000362603846329
$ws.Range('I6').Value = 0.9954146511477043
$ws.Range('J6').Value = 0.9927867998959067
$ws.Range('K6').Value = 0.9919530836898475
$ws.Range('L6').Value = 0.9967543493440122
$ws.Range('M6').Value = 0.9947705746199594
$ws.Range('N6').Value = 0.9947705746199594
$ws.Range('O6').Value = 0.9932495915662796
$ws.Range('P6').Value = 0.9938314109765888
$ws.Range('Q6').Value = 0.9938314109765888
$ws.Range('R6').Value = 0.9933618291549035
$ws.Range('S6').Value = 0.9933618291549035
$ws.Range('T6').Value = 0.9945798522304532
$ws.Range('A7').Value = 5
$ws.Range('B7').Value = 'OffsetF'
$ws.Range('C7').Value = 1.414749273477291
$ws.Range('D7').Value = 0.5717893320643452
$ws.Range('E7').Value = 0.9390975079181278
$ws.Range('F7').Value = 1.414749273477291
$ws.Range('G7').Value = 0.8418046483612815
$ws.Range('H7').Value = 0.8029581348717065
$ws.Range('I7').Value = 1.069788694488562
$ws.Range('J7').Value = 0.5717893320643452
$ws.Range('K7').Value = 1.414749273477291
$ws.Range('L7').Value = 0.9390975079181278
$ws.Range('M7').Value = 0.7554434199912365
$ws.Range('N7').Value = 0.7554434199912365
$ws.Range('O7').Value = 0.784230496114585
$ws.Range('P7').Value = 0.9752120378199214
$ws.Range('Q7').Value = 0.9752120378199214
$ws.Range('R7').Value = 1.085096346734264
$ws.Range('S7').Value = 1.085096346734264
$ws.Range('T7').Value = 0.940031265196886
$ws.Range('A8').Value = 6
$ws.Range('B8').Value = 'OffsetA'
$ws.Range('C8').Value = 0.8516472369506882
$ws.Range('D8').Value = 0.9664389278803046
$ws.Range('E8').Value = 1.077757910552829
$ws.Range('F8').Value = 0.8516472369506882
$ws.Range('G8').Value = 0.9625287716695512
$ws.Range('H8').Value = 1.176392152622551
$ws.Range('I8').Value = 1.009033221104517
$ws.Range('J8').Value = 0.9664389278803046
$ws.Range('K8').Value = 0.8516472369506882
$ws.Range('L8').Value = 1.077757910552829
$ws.Range('M8').Value = 1.022098419216567
$ws.Range('N8').Value = 1.022098419216567
$ws.Range('O8').Value = 1.002241870034228
$ws.Range('P8').Value = 0.9652813584612741
$ws.Range('Q8').Value = 0.9652813584612741
$ws.Range('R8').Value = 0.9368728280836276
$ws.Range('S8').Value = 0.9368728280836276
$ws.Range('T8').Value = 1.007299703463407
$ws.Range('A9').Value = 7
$ws.Range('B9').Value = 'RD Single'
$ws.Range('C9').Value = 1.97
$ws.Range('D9').Value = 0.21
$ws.Range('E9').Value = 0.83
$ws.Range('F9').Value = 1.97
$ws.Range('G9').Value = 0.64
$ws.Range('H9').Value = 0.7
$ws.Range('I9').Value = 1.14
$ws.Range('J9').Value = 0.21
$ws.Range('K9').Value = 1.97
$ws.Range('L9').Value = 0.83
$ws.Range('M9').Value = 0.52
$ws.Range('N9').Value = 0.52
$ws.Range('O9').Value = 0.5600000000000001
$ws.Range('P9').Value = 1.003333333333333
$ws.Range('Q9').Value = 1.003333333333333
$ws.Range('R9').Value = 1.245
$ws.Range('S9').Value = 1.245
$ws.Range('T9').Value = 0.915
$ws.Range('A10').Value = 8
$ws.Range('B10').Value = 'TD Single'
$ws.Range('C10').Value = 0.74
$ws.Range('D10').Value = 0.09
$ws.Range('E10').Value = 1.46
$ws.Range('F10').Value = 0.74
$ws.Range('G10').Value = 0.29
$ws.Range('H10').Value = 2.58
$ws.Range('I10').Value = 1.27
$ws.Range('J10').Value = 0.09
$ws.Range('K10').Value = 0.74
$ws.Range('L10').Value = 1.46
$ws.Range('M10').Value = 0.775
$ws.Range('N10').Value = 0.775
$ws.Range('O10').Value = 0.6133333333333334
$ws.Range('P10').Value = 0.7633333333333333
$ws.Range('Q10').Value = 0.7633333333333333
$ws.Range('R10').Value = 0.7575000000000001
$ws.Range('S10').Value = 0.7575000000000001
$ws.Range('T10').Value = 1.071666666666667
$ws.Range('A11').Value = 9
$ws.Range('B11').Value = 'HexGrid-90degTilt5degRes'
$ws.Range('C11').Value = 0.9953760915961852
$ws.Range('D11').Value = 0.9942690565328817
$ws.Range('E11').Value = 0.9945233956007195
$ws.Range('F11').Value = 0.9953760915961852
$ws.Range('G11').Value = 0.9919540099929121
$ws.Range('H11').Value = 0.9942273845604775
$ws.Range('I11').Value = 0.9946471066514994
$ws.Range('J11').Value = 0.9942690565328817
$ws.Range('K11').Value = 0.9953760915961852
$ws.Range('L11').Value = 0.9945233956007195
$ws.Range('M11').Value = 0.9943962260668006
$ws.Range('N11').Value = 0.9943962260668006
$ws.Range('O11').Value = 0.993582154042171
$ws.Range('P11').Value = 0.9947228479099288
$ws.Range('Q11').Value = 0.9947228479099288
$ws.Range('R11').Value = 0.9948861588314929
$ws.Range('S11').Value = 0.9948861588314929
$ws.Range('T11').Value = 0.9941661741557793

# Apply the same bold/centered/bordered header style already used on row 1 and
# column A to the newly added header cells, by copying formatting (not values)
# from an adjacent cell that already carries that style.
$ws.Range("S1").Copy()
$ws.Range("T1").PasteSpecial(-4122)

$ws.Range("A7").Copy()
$ws.Range("A8:A11").PasteSpecial(-4122)

$excel.CutCopyMode = $false

